# Update countries & provincias Spain
# - refresh "last updated" timestamp
# - refresh COVID case counters for several countries
# - two pairs of countries swapped rank order (their name + stats trade rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "Datos actualizados a 5 de Septiembre de 2020 a las 02:40"
$ws.Range("B4").Value = 6387278
$ws.Range("C4").Value = 52034
$ws.Range("D4").Value = 3630677
$ws.Range("E4").Value = 2564565
$ws.Range("G4").Value = 978
$ws.Range("H4").Value = 192036
$ws.Range("B13").Value = 461882
$ws.Range("C13").Value = 10684
$ws.Range("E13").Value = 120638
$ws.Range("G13").Value = 262
$ws.Range("H13").Value = 9623
$ws.Range("B59").Value = 44777
$ws.Range("C59").Value = 64
$ws.Range("D59").Value = 43693
$ws.Range("E59").Value = 801
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = 283
$ws.Range("A139").Value = "Aruba"
$ws.Range("B139").Value = 2358
$ws.Range("C139").Value = 66
$ws.Range("D139").Value = 1119
$ws.Range("E139").Value = 1226
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 13
$ws.Range("A140").Value = "Jordania"
$ws.Range("B140").Value = 2301
$ws.Range("C140").Value = 68
$ws.Range("D140").Value = 1676
$ws.Range("E140").Value = 609
$ws.Range("G140").Value = 1
$ws.Range("H140").Value = 16
$ws.Range("B141").Value = 2245
$ws.Range("C141").Value = 40
$ws.Range("E141").Value = 1080
$ws.Range("G141").Value = 4
$ws.Range("H141").Value = 38
$ws.Range("A162").Value = "Birmania"
$ws.Range("B162").Value = 1171
$ws.Range("C162").Value = 60
$ws.Range("D162").Value = 359
$ws.Range("E162").Value = 805
$ws.Range("G162").Value = 1
$ws.Range("H162").Value = 7
$ws.Range("A163").Value = "Belice"
$ws.Range("B163").Value = 1152
$ws.Range("C163").Value = 34
$ws.Range("D163").Value = 271
$ws.Range("E163").Value = 866
$ws.Range("G163").Value = 2
$ws.Range("H163").Value = 15
$ws.Range("A164").Value = "Lesoto"
$ws.Range("B164").Value = 1148
$ws.Range("C164").Value = 63
$ws.Range("D164").Value = 528
$ws.Range("E164").Value = 589
$ws.Range("H164").Value = 31
$ws.Range("D167").Value = 859
$ws.Range("E167").Value = 23
$ws.Range("B168").Value = 758
$ws.Range("C168").Value = 4
$ws.Range("E168").Value = 642
$ws.Range("B171").Value = 694
$ws.Range("C171").Value = 72
$ws.Range("D171").Value = 360
$ws.Range("E171").Value = 334
$ws.Range("A173").Value = "San Martin (Parte Holandesa)"
$ws.Range("B173").Value = 511
$ws.Range("C173").Value = 7
$ws.Range("D173").Value = 302
$ws.Range("E173").Value = 190
$ws.Range("H173").Value = 19
$ws.Range("A174").Value = "Tanzania"
$ws.Range("B174").Value = 509
$ws.Range("D174").Value = 183
$ws.Range("E174").Value = 305
$ws.Range("H174").Value = 21
$ws.Range("B177").Value = 451
$ws.Range("C177").Value = 3
$ws.Range("E177").Value = 93
$ws.Range("B189").Value = 178
$ws.Range("C189").Value = 1
$ws.Range("D189").Value = 153
$ws.Range("E189").Value = 18
